# Add CO2 trade links to the islands (Fueltrade / Sheet3)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fueltrade")

# New rows 36-43: B (region1), C (region2), D/E/F (commodity = CO2)
$ws.Range("B36").Value = "DKISLBH"
$ws.Range("C36").Value = "DKE"
$ws.Range("D36").Value = "CO2"
$ws.Range("E36").Value = "CO2"
$ws.Range("F36").Value = "CO2"

$ws.Range("B37").Value = "DKISLBH"
$ws.Range("C37").Value = "DKE"
$ws.Range("D37").Value = "CO2"
$ws.Range("E37").Value = "CO2"
$ws.Range("F37").Value = "CO2"

$ws.Range("B38").Value = "DKISL1"
$ws.Range("C38").Value = "DKW"
$ws.Range("D38").Value = "CO2"
$ws.Range("E38").Value = "CO2"
$ws.Range("F38").Value = "CO2"

$ws.Range("B39").Value = "DKISL1"
$ws.Range("C39").Value = "DKW"
$ws.Range("D39").Value = "CO2"
$ws.Range("E39").Value = "CO2"
$ws.Range("F39").Value = "CO2"

$ws.Range("B40").Value = "DKISL2"
$ws.Range("C40").Value = "DKW"
$ws.Range("D40").Value = "CO2"
$ws.Range("E40").Value = "CO2"
$ws.Range("F40").Value = "CO2"

$ws.Range("B41").Value = "DKISL2"
$ws.Range("C41").Value = "DKW"
$ws.Range("D41").Value = "CO2"
$ws.Range("E41").Value = "CO2"
$ws.Range("F41").Value = "CO2"

$ws.Range("B42").Value = "DKISL3"
$ws.Range("C42").Value = "DKW"
$ws.Range("D42").Value = "CO2"
$ws.Range("E42").Value = "CO2"
$ws.Range("F42").Value = "CO2"

$ws.Range("B43").Value = "DKISL3"
$ws.Range("C43").Value = "DKW"
$ws.Range("D43").Value = "CO2"
$ws.Range("E43").Value = "CO2"
$ws.Range("F43").Value = "CO2"

# Column G (TradeLink names) - entered in this particular order
$ws.Range("G42").Value = "TB_CO2_DKISL3_DKW_01"
$ws.Range("G36").Value = "TB_CO2_DKISLBH_DKE_01"
$ws.Range("G38").Value = "TB_CO2_DKISL1_DKW_01"
$ws.Range("G40").Value = "TB_CO2_DKISL2_DKW_01"
$ws.Range("G37").Value = "TB_CO2_DKISLBH_DKE_02"
$ws.Range("G39").Value = "TB_CO2_DKISL1_DKW_02"
$ws.Range("G41").Value = "TB_CO2_DKISL2_DKW_02"
$ws.Range("G43").Value = "TB_CO2_DKISL3_DKW_02"

# Column width adjustment for column G (from diff: bestFit removed, width set explicitly)
$ws.Columns.Item(7).ColumnWidth = 22.8

# Page setup now explicitly recorded for this sheet (portrait orientation)
$ws.PageSetup.Orientation = 1

# Update the selection to match the final state shown in the diff
$ws.Range("K42").Select() | Out-Null
